$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.450.89"
$ws.Range("E2").Value = "  +3.48%  "

# Row 3
$ws.Range("D3").Value = "3.490.60"
$ws.Range("E3").Value = "  +2.66%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'589.69"
$ws.Range("E5").Value = "  +2.65%  "

# Row 6
$ws.Range("D6").Value = "'168.00"
$ws.Range("E6").Value = "  +3.14%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("D8").Value = "3.487.53"
$ws.Range("E8").Value = "  +2.49%  "

# Row 9
$ws.Range("D9").Value = "'0.590"
$ws.Range("E9").Value = "  +7.33%  "

# Row 10
$ws.Range("D10").Value = "'7.33"
$ws.Range("E10").Value = "  +0.41%  "

# Row 11
$ws.Range("E11").Value = "  +5.75%  "

# Row 12
$ws.Range("D12").Value = "'0.430"
$ws.Range("E12").Value = "  +2.38%  "

# Row 13
$ws.Range("D13").Value = "4.093.31"
$ws.Range("E13").Value = "  +2.65%  "

# Row 14
$ws.Range("E14").Value = "  -0.55%  "

# Row 15
$ws.Range("D15").Value = "'27.94"
$ws.Range("E15").Value = "  +4.13%  "

# Row 16
$ws.Range("D16").Value = "66.473.20"
$ws.Range("E16").Value = "  +3.41%  "

# Row 17
$ws.Range("D17").Value = "'0.0000176"
$ws.Range("E17").Value = "  +2.32%  "

# Row 18
$ws.Range("D18").Value = "3.485.44"
$ws.Range("E18").Value = "  +2.51%  "

# Row 19
$ws.Range("D19").Value = "'6.25"
$ws.Range("E19").Value = "  +1.96%  "

# Row 20
$ws.Range("E20").Value = "  +3.51%  "

# Row 21
$ws.Range("D21").Value = "'389.41"
$ws.Range("E21").Value = "  +4.07%  "

# Row 22
$ws.Range("D22").Value = "'7.87"
$ws.Range("E22").Value = "  +1.16%  "

# Row 23
$ws.Range("D23").Value = "'72.63"
$ws.Range("E23").Value = "  +3.10%  "

# Row 24
$ws.Range("D24").Value = "'0.998"
$ws.Range("E24").Value = "  -0.32%  "

# Row 25
$ws.Range("E25").Value = "  +3.36%  "

# Row 26
$ws.Range("D26").Value = "'0.0000121"
$ws.Range("E26").Value = "  +5.55%  "

# Row 27
$ws.Range("D27").Value = "'10.16"
$ws.Range("E27").Value = "  +7.45%  "

# Row 28
$ws.Range("E28").Value = "  +0.25%  "

# Row 29
$ws.Range("E29").Value = "  -0.41%  "

# Row 30
$ws.Range("D30").Value = "'6.30"
$ws.Range("E30").Value = "  +4.09%  "

# Row 31
$ws.Range("E31").Value = "  +3.12%  "

# Row 32
$ws.Range("D32").Value = "'2.04"
$ws.Range("E32").Value = "  +2.14%  "

# Row 33
$ws.Range("D33").Value = "'23.57"
$ws.Range("E33").Value = "  +3.43%  "

# Row 34
$ws.Range("D34").Value = "'7.29"
$ws.Range("E34").Value = "  +4.25%  "

# Row 35
$ws.Range("D35").Value = "'1.58"
$ws.Range("E35").Value = "  +7.17%  "

# Row 36
$ws.Range("D36").Value = "'162.65"
$ws.Range("E36").Value = "  +2.30%  "

# Row 37
$ws.Range("D37").Value = "'0.898"
$ws.Range("E37").Value = "  +4.24%  "

# Row 38
$ws.Range("E38").Value = "  +4.62%  "

# Row 39
$ws.Range("D39").Value = "'6.78"
$ws.Range("E39").Value = "  +5.07%  "

# Row 40
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'4.61"
$ws.Range("E40").Value = "  +6.14%  "

# Row 41
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.0736"
$ws.Range("E41").Value = "  +2.52%  "

# Row 42
$ws.Range("D42").Value = "'26.24"
$ws.Range("E42").Value = "  +2.37%  "

# Row 43
$ws.Range("D43").Value = "2.779.64"
$ws.Range("E43").Value = "  +0.54%  "

# Row 44
$ws.Range("D44").Value = "'26.45"
$ws.Range("E44").Value = "  +2.16%  "

# Row 45
$ws.Range("D45").Value = "'42.65"
$ws.Range("E45").Value = "  +0.31%  "

# Row 46
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "'2.50"
$ws.Range("E46").Value = "  +3.23%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0308"
$ws.Range("E47").Value = "  +1.88%  "

# Row 48
$ws.Range("D48").Value = "'343.03"
$ws.Range("E48").Value = "  +4.40%  "

# Row 49
$ws.Range("D49").Value = "'1.07"
$ws.Range("E49").Value = "  +2.79%  "

# Row 50
$ws.Range("D50").Value = "'33.11"
$ws.Range("E50").Value = "  +10.35%  "

# Row 51
$ws.Range("D51").Value = "'0.855"
$ws.Range("E51").Value = "  +5.42%  "
